$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New log entry: "Dia 2" / night study session ---
# Rows 9-10 already exist as an empty templated pair (C9:C10, E9:E10, F9:F10 merged).
# D9:D10 is not yet merged because that pair never had a "Tema" entry filled in.

# 1) Merge D9:D10 first so the subsequent format paste can settle onto the final
#    merged range (mirrors how the other data rows - 5:6, 7:8, 11:12, ... - are built).
$ws.Range("D9:D10").Merge()

# 2) Copy the visual formatting (fonts/borders/alignment) from the previous entry
#    (rows 7-8) down into the new rows 9-10, without touching the text we are about
#    to add, and without disturbing the merges already in place.
$ws.Range("C7:F8").Copy()
$ws.Range("C9:F10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# 3) Fill in the values for the new entry.
$ws.Range("C9").Value = $ws.Range("C7").Value2
$ws.Range("D9").Value = "Palabras claves cons final dynamic`nManipulación de String"
$ws.Range("E9").Value = $ws.Range("E7").Value2
$ws.Range("F9").Value = $ws.Range("F7").Value2

# 4) Match row 9's height to the standard (unwrapped) entry-row height.
$ws.Range("A9").RowHeight = 14.4
